$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 22.24854560026249
$ws.Range("C2").Value = 8.108704586118678
$ws.Range("D2").Value = 13.57497534017594
$ws.Range("E2").Value = 13.8832556097768
$ws.Range("G2").Value = 3.713878746361673
$ws.Range("J2").Value = 8.436695558965017
$ws.Range("L2").Value = 12.73779624624096
$ws.Range("O2").Value = 32.73553132502259

$ws.Range("B3").Value = 21.82572601128805
$ws.Range("C3").Value = 7.792191563136518
$ws.Range("D3").Value = 13.58785991460878
$ws.Range("E3").Value = 13.92413729900321
$ws.Range("G3").Value = 3.716744758009138
$ws.Range("J3").Value = 8.447732265160848
$ws.Range("L3").Value = 12.72328869931621
$ws.Range("O3").Value = 32.8468727494686

$ws.Range("B4").Value = 21.56548130754456
$ws.Range("C4").Value = 7.590098063782095
$ws.Range("D4").Value = 13.59831687860034
$ws.Range("E4").Value = 13.9511021607638
$ws.Range("G4").Value = 3.718597094492405
$ws.Range("J4").Value = 8.45486676981751
$ws.Range("L4").Value = 12.71616038669058
$ws.Range("O4").Value = 32.92335162319041

$ws.Range("B5").Value = 21.45940925709928
$ws.Range("C5").Value = 7.505877518332971
$ws.Range("D5").Value = 13.60321806014536
$ws.Range("E5").Value = 13.96255944193575
$ws.Range("G5").Value = 3.719375300108258
$ws.Range("J5").Value = 8.457864430246151
$ws.Range("L5").Value = 12.71370528815604
$ws.Range("O5").Value = 32.95655050609444

$ws.Range("B6").Value = 21.44179895766273
$ws.Range("C6").Value = 7.491782731729628
$ws.Range("D6").Value = 13.60407053331834
$ws.Range("E6").Value = 13.96449024378046
$ws.Range("G6").Value = 3.719505934005481
$ws.Range("J6").Value = 8.458367651392219
$ws.Range("L6").Value = 12.71332484751771
$ws.Range("O6").Value = 32.9621857297922

$ws.Range("B7").Value = 21.56405067947683
$ws.Range("C7").Value = 7.58896967197977
$ws.Range("D7").Value = 13.59838038728433
$ws.Range("E7").Value = 13.95125477874536
$ws.Range("G7").Value = 3.718607494942537
$ws.Range("J7").Value = 8.454906831343614
$ws.Range("L7").Value = 12.71612545248586
$ws.Range("O7").Value = 32.92379113294889

$ws.Range("B8").Value = 22.10297077896798
$ws.Range("C8").Value = 8.001227093736809
$ws.Range("D8").Value = 13.57888945185196
$ws.Range("E8").Value = 13.89696507840411
$ws.Range("G8").Value = 3.714847777594849
$ws.Range("J8").Value = 8.44042690296021
$ws.Range("L8").Value = 12.73242618069929
$ws.Range("O8").Value = 32.77223323641574

$ws.Range("B9").Value = 23.14885806716237
$ws.Range("C9").Value = 8.74498346334263
$ws.Range("D9").Value = 13.56087345282418
$ws.Range("E9").Value = 13.80527314458794
$ws.Range("G9").Value = 3.708206022766634
$ws.Range("J9").Value = 8.414858545632026
$ws.Range("L9").Value = 12.77840129428469
$ws.Range("O9").Value = 32.53973817548081

$ws.Range("B10").Value = 23.90298536234462
$ws.Range("C10").Value = 9.248416801970839
$ws.Range("D10").Value = 13.55995398585839
$ws.Range("E10").Value = 13.74689100345928
$ws.Range("G10").Value = 3.703766877473886
$ws.Range("J10").Value = 8.397778211658604
$ws.Range("L10").Value = 12.82056208243096
$ws.Range("O10").Value = 32.40879981652022

$ws.Range("B11").Value = 24.24144816882917
$ws.Range("C11").Value = 9.467514117155149
$ws.Range("D11").Value = 13.5622058571662
$ws.Range("E11").Value = 13.72227806303117
$ws.Range("G11").Value = 3.701841968171236
$ws.Range("J11").Value = 8.390374150974719
$ws.Range("L11").Value = 12.84152423397441
$ws.Range("O11").Value = 32.35797721202609

$ws.Range("B12").Value = 24.36883736361554
$ws.Range("C12").Value = 9.549014613693023
$ws.Range("D12").Value = 13.56344176766052
$ws.Range("E12").Value = 13.71323724270317
$ws.Range("G12").Value = 3.701126558196334
$ws.Range("J12").Value = 8.387622738380482
$ws.Range("L12").Value = 12.84971480663374
$ws.Range("O12").Value = 32.33999567863214

$ws.Range("B13").Value = 24.34143830087243
$ws.Range("C13").Value = 9.53152784959126
$ws.Range("D13").Value = 13.5631585653018
$ws.Range("E13").Value = 13.71517191412596
$ws.Range("G13").Value = 3.701280034796529
$ws.Range("J13").Value = 8.388212980506673
$ws.Range("L13").Value = 12.84793964487177
$ws.Range("O13").Value = 32.34381200722174

$ws.Range("B14").Value = 24.2519448296015
$ws.Range("C14").Value = 9.474248823040565
$ws.Range("D14").Value = 13.56229986215981
$ws.Range("E14").Value = 13.72152866738298
$ws.Range("G14").Value = 3.701782840599784
$ws.Range("J14").Value = 8.390146743065454
$ws.Range("L14").Value = 12.84219303577671
$ws.Range("O14").Value = 32.35647249255868

$ws.Range("B15").Value = 24.19702259377595
$ws.Range("C15").Value = 9.438971624733536
$ws.Range("D15").Value = 13.561823755291
$ws.Range("E15").Value = 13.72545876485582
$ws.Range("G15").Value = 3.702092581045907
$ws.Range("J15").Value = 8.391338037376581
$ws.Range("L15").Value = 12.83870586114609
$ws.Range("O15").Value = 32.36439219660922

$ws.Range("B16").Value = 23.88076328095113
$ws.Range("C16").Value = 9.233895086480686
$ws.Range("D16").Value = 13.55986048864248
$ws.Range("E16").Value = 13.74853864779947
$ws.Range("G16").Value = 3.703894569606884
$ws.Range("J16").Value = 8.398269423785438
$ws.Range("L16").Value = 12.81922772103079
$ws.Range("O16").Value = 32.41229769250517

$ws.Range("B17").Value = 23.68548528596201
$ws.Range("C17").Value = 9.105516207091325
$ws.Range("D17").Value = 13.55933948611367
$ws.Range("E17").Value = 13.76319554022485
$ws.Range("G17").Value = 3.705024176537964
$ws.Range("J17").Value = 8.402615122045249
$ws.Range("L17").Value = 12.80773259664993
$ws.Range("O17").Value = 32.44393012621878

$ws.Range("B18").Value = 23.57274087408815
$ws.Range("C18").Value = 9.03074464922998
$ws.Range("D18").Value = 13.55929120684372
$ws.Range("E18").Value = 13.77180891652562
$ws.Range("G18").Value = 3.705682793922618
$ws.Range("J18").Value = 8.405149105887375
$ws.Range("L18").Value = 12.80128896686908
$ws.Range("O18").Value = 32.46294670894208

$ws.Range("B19").Value = 23.53449817329658
$ws.Range("C19").Value = 9.005269659666737
$ws.Range("D19").Value = 13.55931805670645
$ws.Range("E19").Value = 13.77475671422446
$ws.Range("G19").Value = 3.705907320867727
$ws.Range("J19").Value = 8.406012994807995
$ws.Range("L19").Value = 12.79913623938448
$ws.Range("O19").Value = 32.46952643972642

$ws.Range("B20").Value = 23.70631789946866
$ws.Range("C20").Value = 9.119279072267117
$ws.Range("D20").Value = 13.55936893504511
$ws.Range("E20").Value = 13.7616163381369
$ws.Range("G20").Value = 3.704903007650876
$ws.Range("J20").Value = 8.402148951213647
$ws.Range("L20").Value = 12.80893890483824
$ws.Range("O20").Value = 32.44047762936766

$ws.Range("B21").Value = 24.27825329506388
$ws.Range("C21").Value = 9.491113167908626
$ws.Range("D21").Value = 13.56254169219447
$ws.Range("E21").Value = 13.7196539489124
$ws.Range("G21").Value = 3.701634788245128
$ws.Range("J21").Value = 8.389577331630534
$ws.Range("L21").Value = 12.8438741278605
$ws.Range("O21").Value = 32.35271944710707

$ws.Range("B22").Value = 24.64745338456574
$ws.Range("C22").Value = 9.725563268323704
$ws.Range("D22").Value = 13.56684819006281
$ws.Range("E22").Value = 13.69385854798161
$ws.Range("G22").Value = 3.699577538370885
$ws.Range("J22").Value = 8.381666032763395
$ws.Range("L22").Value = 12.86817713100626
$ws.Range("O22").Value = 32.30273419827654

$ws.Range("B23").Value = 24.45086206184478
$ws.Range("C23").Value = 9.601228468854845
$ws.Range("D23").Value = 13.56434573700846
$ws.Range("E23").Value = 13.70747700193255
$ws.Range("G23").Value = 3.700668352993447
$ws.Range("J23").Value = 8.385860623769888
$ws.Range("L23").Value = 12.85507287064633
$ws.Range("O23").Value = 32.32873576947995

$ws.Range("B24").Value = 23.69690094992506
$ws.Range("C24").Value = 9.113059883068114
$ws.Range("D24").Value = 13.55935483850576
$ws.Range("E24").Value = 13.76232971309652
$ws.Range("G24").Value = 3.704957759467206
$ws.Range("J24").Value = 8.402359596180276
$ws.Range("L24").Value = 12.80839301845046
$ws.Range("O24").Value = 32.44203591608379

$ws.Range("B25").Value = 22.86789500993197
$ws.Range("C25").Value = 8.5511045113445
$ws.Range("D25").Value = 13.56358338223451
$ws.Range("E25").Value = 13.82849913119362
$ws.Range("G25").Value = 3.709925060508782
$ws.Range("J25").Value = 8.421474763381678
$ws.Range("L25").Value = 12.76447935898421
$ws.Range("O25").Value = 32.5956614732964
